$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A5").Value = "노래한장앱"
$ws.Range("B5").Value = "시니언노래동영상앱"
$ws.Range("C5").Value = "image4.jpg"
